# Update the Cck-Cckar LR-pairs sheet with refreshed TPM-derived values.
# The "Target cluster" set gained a new cluster ("Resolving-Mac"), so the
# NATMI output table grew from 4 target clusters x 3 sending clusters
# (12 data rows) to 5 target clusters x 3 sending clusters (15 data rows),
# and every numeric column was recomputed against the new TPM values.
#
# Columns (A..T):
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E..J ligand-side stats, K..P receptor-side stats, Q..T edge-weight stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 is the header; rows 2-16 are the 3 (sending cluster) x 5 (target cluster) combinations.
$data = @(
    @('Sending cluster', 'Ligand symbol', 'Receptor symbol', 'Target cluster', 'Ligand-expressing cells', 'Ligand detection rate', 'Ligand average expression value', 'Ligand total expression value', 'Ligand derived specificity of average expression value', 'Ligand derived specificity of total expression value', 'Receptor-expressing cells', 'Receptor detection rate', 'Receptor average expression value', 'Receptor total expression value', 'Receptor derived specificity of average expression value', 'Receptor derived specificity of total expression value', 'Edge average expression weight', 'Edge total expression weight', 'Edge average expression derived specificity', 'Edge total expression derived specificity'),
    @('ECs', 'Cck', 'Cckar', 'FAPs', 1, 0.3333333333333333, 0.234784, 0.704352, 0.1954244153445811, 0.1954244153445811, 3, 1, 0.5512466666666667, 1.65374, 0.5047325351567603, 0.5047325351567603, 0.1294238973866667, 1.16481507648, 0.09863706058839809, 0.09863706058839809),
    @('ECs', 'Cck', 'Cckar', 'Inflammatory-Mac', 1, 0.3333333333333333, 0.234784, 0.704352, 0.1954244153445811, 0.1954244153445811, 2, 0.6666666666666666, 0.2660283333333333, 0.7980849999999999, 0.2435808925953191, 0.2435808925953191, 0.06245919621333332, 0.5621327659199999, 0.04760165352455143, 0.04760165352455143),
    @('ECs', 'Cck', 'Cckar', 'MuSCs', 1, 0.3333333333333333, 0.234784, 0.704352, 0.1954244153445811, 0.1954244153445811, 1, 0.3333333333333333, 0.018172, 0.054516, 0.01663864869121261, 0.01663864869121261, 0.004266494848, 0.038398453632, 0.003251598192604104, 0.003251598192604103),
    @('ECs', 'Cck', 'Cckar', 'Neutrophils', 1, 0.3333333333333333, 0.234784, 0.704352, 0.1954244153445811, 0.1954244153445811, 1, 0.3333333333333333, 0.1763026666666667, 0.528908, 0.1614262675539636, 0.1614262675539636, 0.04139304529066667, 0.372537407616, 0.03154663395799126, 0.03154663395799125),
    @('ECs', 'Cck', 'Cckar', 'Resolving-Mac', 1, 0.3333333333333333, 0.234784, 0.704352, 0.1954244153445811, 0.1954244153445811, 1, 0.3333333333333333, 0.08040633333333333, 0.241219, 0.07362165600274442, 0.07362165600274442, 0.01887812056533333, 0.169903085088, 0.0143874690810362, 0.0143874690810362),
    @('FAPs', 'Cck', 'Cckar', 'FAPs', 2, 0.6666666666666666, 0.6129226666666666, 1.838768, 0.5101712799201602, 0.5101712799201601, 3, 1, 0.5512466666666667, 1.65374, 0.5047325351567603, 0.5047325351567603, 0.3378715769244444, 3.04084419232, 0.2575000434782717, 0.2575000434782716),
    @('FAPs', 'Cck', 'Cckar', 'Inflammatory-Mac', 2, 0.6666666666666666, 0.6129226666666666, 1.838768, 0.5101712799201602, 0.5101712799201601, 2, 0.6666666666666666, 0.2660283333333333, 0.7980849999999999, 0.2435808925953191, 0.2435808925953191, 0.1630547954755555, 1.46749315928, 0.124267975739449, 0.124267975739449),
    @('FAPs', 'Cck', 'Cckar', 'MuSCs', 2, 0.6666666666666666, 0.6129226666666666, 1.838768, 0.5101712799201602, 0.5101712799201601, 1, 0.3333333333333333, 0.018172, 0.054516, 0.01663864869121261, 0.01663864869121261, 0.01113803069866667, 0.100242276288, 0.008488560698937837, 0.008488560698937834),
    @('FAPs', 'Cck', 'Cckar', 'Neutrophils', 2, 0.6666666666666666, 0.6129226666666666, 1.838768, 0.5101712799201602, 0.5101712799201601, 1, 0.3333333333333333, 0.1763026666666667, 0.528908, 0.1614262675539636, 0.1614262675539636, 0.1080599005937778, 0.972539105344, 0.08235504553073986, 0.08235504553073983),
    @('FAPs', 'Cck', 'Cckar', 'Resolving-Mac', 2, 0.6666666666666666, 0.6129226666666666, 1.838768, 0.5101712799201602, 0.5101712799201601, 1, 0.3333333333333333, 0.08040633333333333, 0.241219, 0.07362165600274442, 0.07362165600274442, 0.04928286424355555, 0.443545778192, 0.03755965447276187, 0.03755965447276186),
    @('MuSCs', 'Cck', 'Cckar', 'FAPs', 3, 1, 0.353699, 1.061097, 0.2944043047352588, 0.2944043047352587, 3, 1, 0.5512466666666667, 1.65374, 0.5047325351567603, 0.5047325351567603, 0.1949753947533334, 1.75477855278, 0.1485954310900906, 0.1485954310900906),
    @('MuSCs', 'Cck', 'Cckar', 'Inflammatory-Mac', 3, 1, 0.353699, 1.061097, 0.2944043047352588, 0.2944043047352587, 2, 0.6666666666666666, 0.2660283333333333, 0.7980849999999999, 0.2435808925953191, 0.2435808925953191, 0.09409395547166667, 0.846845599245, 0.07171126333131866, 0.07171126333131865),
    @('MuSCs', 'Cck', 'Cckar', 'MuSCs', 3, 1, 0.353699, 1.061097, 0.2944043047352588, 0.2944043047352587, 1, 0.3333333333333333, 0.018172, 0.054516, 0.01663864869121261, 0.01663864869121261, 0.006427418228000001, 0.05784676405200001, 0.004898489799670673, 0.004898489799670671),
    @('MuSCs', 'Cck', 'Cckar', 'Neutrophils', 3, 1, 0.353699, 1.061097, 0.2944043047352588, 0.2944043047352587, 1, 0.3333333333333333, 0.1763026666666667, 0.528908, 0.1614262675539636, 0.1614262675539636, 0.06235807689733335, 0.5612226920760002, 0.04752458806523253, 0.04752458806523251),
    @('MuSCs', 'Cck', 'Cckar', 'Resolving-Mac', 3, 1, 0.353699, 1.061097, 0.2944043047352588, 0.2944043047352587, 1, 0.3333333333333333, 0.08040633333333333, 0.241219, 0.07362165600274442, 0.07362165600274442, 0.02843963969366667, 0.255956757243, 0.02167453244894637, 0.02167453244894636)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 1
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $data[$r][$c]
    }
}
